$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old used range completely (A1:E12, plus the new column F) so
# stale cells (old row 12, old headers/values) don't linger after the
# layout change.
$ws.Range("A1:F12").ClearContents()

# --- Row 1: section headers ---
$ws.Cells.Item(1,1).Value = "物品编号"
$ws.Cells.Item(1,3).Value = "怪物编号"
$ws.Cells.Item(1,5).Value = "事件编号"

# --- Column A/B: item list (rows 2-11) ---
$itemNames = @("黄钥匙","蓝钥匙","红钥匙","红血瓶","蓝血瓶","红宝石","蓝宝石","怪物手册","楼层传送器","圣水")
for ($i = 0; $i -lt $itemNames.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $itemNames[$i]
}

# --- Column C/D: monster list (rows 2-4) ---
$monsterNames = @("绿史莱姆","红史莱姆","黑史莱姆")
for ($i = 0; $i -lt $monsterNames.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $i + 1
    $ws.Cells.Item($r, 4).Value = $monsterNames[$i]
}

# --- Column E/F: event list (rows 2-5) ---
$eventNames = @("黄门","蓝门","红门","铁门")
for ($i = 0; $i -lt $eventNames.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $i + 1
    $ws.Cells.Item($r, 6).Value = $eventNames[$i]
}

$ws.Range("A12").Select()
